$wb = $excel.ActiveWorkbook
$wsGames = $wb.Worksheets.Item("Games")
$wsNext  = $wb.Worksheets.Item("Next")

# --- 1) "Next" sheet: the first upcoming game (45306 @ IND) has now been
#        played, so remove it from the "Next" list. Deleting the row
#        shifts every remaining row up by one, matching the diff.
$wsNext.Rows.Item(2).Delete()

# --- 2) "Games" sheet: append the completed game as the new last row
#        (Game 42), using the stats/result for that game.
$lastRow = $wsGames.UsedRange.Rows.Count
$newRow = $lastRow + 1

$dateFmt = $wsGames.Cells.Item($lastRow, 2).NumberFormat

$wsGames.Cells.Item($newRow, 1).Value2  = 42
$wsGames.Cells.Item($newRow, 2).Value2  = 45306
$wsGames.Cells.Item($newRow, 2).NumberFormat = $dateFmt
$wsGames.Cells.Item($newRow, 3).Value2  = 6
$wsGames.Cells.Item($newRow, 4).Value2  = 132
$wsGames.Cells.Item($newRow, 5).Value2  = 97.8
$wsGames.Cells.Item($newRow, 6).Value2  = 0.623
$wsGames.Cells.Item($newRow, 7).Value2  = 12.7
$wsGames.Cells.Item($newRow, 8).Value2  = 28.6
$wsGames.Cells.Item($newRow, 9).Value2  = 0.383
$wsGames.Cells.Item($newRow, 10).Value2 = 134.9
$wsGames.Cells.Item($newRow, 11).Value2 = "IND"
$wsGames.Cells.Item($newRow, 12).Value2 = 105
$wsGames.Cells.Item($newRow, 13).Value2 = 0.447
$wsGames.Cells.Item($newRow, 14).Value2 = 10.9
$wsGames.Cells.Item($newRow, 15).Value2 = 33.3
$wsGames.Cells.Item($newRow, 16).Value2 = 0.211
$wsGames.Cells.Item($newRow, 17).Value2 = 107.3
$wsGames.Cells.Item($newRow, 18).Value2 = 1
$wsGames.Cells.Item($newRow, 19).Value2 = 1
